# Insert a new daily price record for "Albahaca" just before the existing
# row 190, shifting the old rows 190-208 down to 191-209 (the sheet grows
# from A1:R208 to A1:R209).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 190..208 down by one, opening a blank row 190.
$ws.Rows(190).Insert()

# Populate the newly opened row 190 with the new record's data.
$ws.Range("A190").Value = 9
$ws.Range("B190").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C190").Value = "Metropolitana"
$ws.Range("D190").Value = 44491
$ws.Range("E190").Value = 13
$ws.Range("F190").Value = 100112052
$ws.Range("G190").Value = "Albahaca"
$ws.Range("H190").Value = "Sin especificar"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 97
$ws.Range("K190").Value = 6000
$ws.Range("L190").Value = 7000
$ws.Range("M190").Value = 6495
$ws.Range("N190").Value = "`$/docena de matas"
$ws.Range("O190").Value = "Provincia de Chacabuco"
$ws.Range("P190").Value = 1082
$ws.Range("Q190").Value = 6
$ws.Range("R190").Value = "Hortaliza"
